$d = $word.ActiveDocument

# Locate anchor paragraph: the FAILURE paragraph ending in 1525348404493,
# which is the second-to-last paragraph in the original document.
$anchorIndex = $d.Paragraphs.Count - 1

$items = @(
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Connection details saved successfully"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'Yes' , Header at value is 1 and At Line value is 2"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'No' , Header template value is Blank and At Line value is Blank"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Regular Expression is selected so At Line text box  is not visible and  regular expression text box is visible"; Break = $true },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: At Line is selected so At Line text box  is visible and  regular expression text box is not visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Path /invalid/path/ is invalid."; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Incorrect header template"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Please enter record start line expression"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_FAILURE : No record found in /hadoop/softwares/momosaic_6_7/tmp/uploadpath/1525421291253/CSV.csv for specified criterion"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Connection details saved successfully"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'Yes' , Header at value is 1 and At Line value is 2"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'No' , Header template value is Blank and At Line value is Blank"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Regular Expression is selected so At Line text box  is not visible and  regular expression text box is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: At Line is selected so At Line text box  is visible and  regular expression text box is not visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Path /invalid/path/ is invalid."; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Incorrect header template"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Please enter record start line expression"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_FAILURE : No record found in /hadoop/softwares/momosaic_6_7/tmp/uploadpath/1525672392956/CSV.csv for specified criterion"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Connection details saved successfully"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'Yes' , Header at value is 1 and At Line value is 2"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'No' , Header template value is Blank and At Line value is Blank"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Regular Expression is selected so At Line text box  is not visible and  regular expression text box is visible"; Break = $true },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: At Line is selected so At Line text box  is visible and  regular expression text box is not visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Path /invalid/path/ is invalid."; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Incorrect header template"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Please enter record start line expression"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_FAILURE : No record found in /hadoop/softwares/momosaic_6_7/tmp/uploadpath/1525685358553/CSV.csv for specified criterion"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Connection details saved successfully"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'Yes' , Header at value is 1 and At Line value is 2"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'No' , Header template value is Blank and At Line value is Blank"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Regular Expression is selected so At Line text box  is not visible and  regular expression text box is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: At Line is selected so At Line text box  is visible and  regular expression text box is not visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Path /invalid/path/ is invalid."; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Incorrect header template"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Please enter record start line expression"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_FAILURE : No record found in /hadoop/softwares/momosaic_6_7/tmp/uploadpath/1525696034255/CSV.csv for specified criterion"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Connection details saved successfully"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'Yes' , Header at value is 1 and At Line value is 2"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Header at and At Line is visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: For contains header 'No' , Header template value is Blank and At Line value is Blank"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: Regular Expression is selected so At Line text box  is not visible and  regular expression text box is visible"; Break = $true },
    @{ Text = "D_Data_TCNo_10 to 16_Configuration_Connection:SUCESS: At Line is selected so At Line text box  is visible and  regular expression text box is not visible"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Path /invalid/path/ is invalid."; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Incorrect header template"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16 _SUCESS: Please enter record start line expression"; Break = $false },
    @{ Text = "D_Data_TCNo_10 to 16_FAILURE : No record found in /hadoop/softwares/momosaic_6_7/tmp/uploadpath/1525699444697/CSV.csv for specified criterion"; Break = $false }
)

$r = $d.Paragraphs.Item($anchorIndex).Range
$idx = $anchorIndex
foreach ($item in $items) {
    $r.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs.Item($idx)
    if ($item.Break) {
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>' + $item.Text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
    } else {
        $p.Range.Text = $item.Text
    }
    $r = $p.Range
}

Write-Output "Inserted $($items.Count) paragraphs. Total paragraphs now: $($d.Paragraphs.Count)"
